$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the #ALUTs values for both data rows (column G) from 2129 to 2140.
# Dependent formulas in column I (#ALUTs + #Registers) and column M
# (Area*Delay) recalculate automatically.
$ws.Range("G2").Value = 2140
$ws.Range("G3").Value = 2140

# Update the active selection to match the authored state.
$ws.Range("I6").Select()
